$wb = $excel.ActiveWorkbook

# --- Rename worksheets ---
$wsBenthic = $wb.Worksheets.Item("benthic_algae")
$wsBenthic.Name = "benthic"
$wsInverts = $wb.Worksheets.Item("benthic_inverts")
$wsInverts.Name = "inverts"

# --- Add new invertebrate survey rows to "inverts" sheet ---
# (populate the new "Conch" string reference before the "Macroalgal height"
#  header rename below, so new shared-string entries land in the same order
#  as the target workbook: Conch, then Macroalgal height)
$wsInverts.Range("A2:F2").Copy($wsInverts.Range("A3:F3"))
$wsInverts.Range("A2:F2").Copy($wsInverts.Range("A4:F4"))
$wsInverts.Range("A2:F2").Copy($wsInverts.Range("A5:F5"))
$wsInverts.Range("A2:F2").Copy($wsInverts.Range("A6:F6"))
$wsInverts.Range("A2:F2").Copy($wsInverts.Range("A7:F7"))

$wsInverts.Range("E3").Value = "Conch"
$wsInverts.Range("F3").Value = 1

$wsInverts.Range("E4").Value = "Lobster"
$wsInverts.Range("F4").Value = 2

$wsInverts.Range("C5").Value = 44477
$wsInverts.Range("E5").Value = "Lobster"
$wsInverts.Range("F5").Value = 1

$wsInverts.Range("C6").Value = 44477
$wsInverts.Range("E6").Value = "Conch"
$wsInverts.Range("F6").Value = 2

$wsInverts.Range("C7").Value = 44477
$wsInverts.Range("E7").Value = "Lobster"
$wsInverts.Range("F7").Value = 3

$wsInverts.Range("F7").Select()

# --- Fix header typo on "benthic" sheet: Macroalgae height -> Macroalgal height ---
$wsBenthic.Range("M1").Value = "Macroalgal height"

# --- Make "benthic" the active/selected sheet with its target view state ---
$wsBenthic.Activate()
$wsBenthic.Range("J12").Select()
